$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44701
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 7500
$ws.Range("M2").Value = 7250
$ws.Range("N2").Value = '$/paquete 36 unidades'
$ws.Range("O2").Value = 'Región Metropolitana'
$ws.Range("P2").Value = 201
$ws.Range("Q2").Value = 36

$ws.Range("D3").Value = 44225
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 3400
$ws.Range("L3").Value = 3700
$ws.Range("M3").Value = 3550
$ws.Range("N3").Value = '$/paquete 2 kilos'
$ws.Range("O3").Value = 'Provincia de Diguillín'
$ws.Range("P3").Value = 1775
$ws.Range("Q3").Value = 2

$ws.Range("D5").Value = 44662
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 8500
$ws.Range("M5").Value = 8250
$ws.Range("N5").Value = '$/paquete 36 unidades'
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 229
$ws.Range("Q5").Value = 36

$ws.Range("D7").Value = 44223
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 3500
$ws.Range("L7").Value = 3800
$ws.Range("M7").Value = 3688
$ws.Range("N7").Value = '$/paquete 2 kilos'
$ws.Range("O7").Value = 'Provincia de Diguillín'
$ws.Range("P7").Value = 1844
$ws.Range("Q7").Value = 2

$ws.Range("D8").Value = 44208
$ws.Range("J8").Value = 85
$ws.Range("K8").Value = 3700
$ws.Range("L8").Value = 4000
$ws.Range("M8").Value = 3824
$ws.Range("N8").Value = '$/paquete 2 kilos'
$ws.Range("O8").Value = 'Provincia de Diguillín'
$ws.Range("P8").Value = 1912
$ws.Range("Q8").Value = 2

$ws.Range("D9").Value = 44166
$ws.Range("J9").Value = 70
$ws.Range("K9").Value = 3500
$ws.Range("L9").Value = 4000
$ws.Range("M9").Value = 3679
$ws.Range("N9").Value = '$/paquete 36 unidades'
$ws.Range("O9").Value = 'Región Metropolitana'
$ws.Range("P9").Value = 102
$ws.Range("Q9").Value = 36

$ws.Range("D10").Value = 44161
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 2800
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 2900
$ws.Range("N10").Value = '$/paquete 2 kilos'
$ws.Range("O10").Value = 'Provincia de Diguillín'
$ws.Range("P10").Value = 1450
$ws.Range("Q10").Value = 2

$ws.Range("D11").Value = 44160
$ws.Range("J11").Value = 43
$ws.Range("K11").Value = 3500
$ws.Range("L11").Value = 4000
$ws.Range("M11").Value = 3709
$ws.Range("N11").Value = '$/paquete 36 unidades'
$ws.Range("O11").Value = 'Región Metropolitana'
$ws.Range("P11").Value = 103
$ws.Range("Q11").Value = 36

$ws.Range("D12").Value = 44664
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 8000
$ws.Range("L12").Value = 8500
$ws.Range("M12").Value = 8250
$ws.Range("N12").Value = '$/paquete 36 unidades'
$ws.Range("O12").Value = 'Región Metropolitana'
$ws.Range("P12").Value = 229
$ws.Range("Q12").Value = 36

$ws.Range("D13").Value = 44209
$ws.Range("J13").Value = 150
$ws.Range("K13").Value = 3500
$ws.Range("L13").Value = 4000
$ws.Range("M13").Value = 3767
$ws.Range("N13").Value = '$/paquete 2 kilos'
$ws.Range("O13").Value = 'Provincia de Diguillín'
$ws.Range("P13").Value = 1884
$ws.Range("Q13").Value = 2
